# Update the cryptocurrency price / volume(1h) figures with the latest
# scraped values (GitHub Actions refresh).
#
# Note: a handful of the new "Price" values happen to look like plain
# decimal numbers (e.g. "582.75"). The source data stores every Price /
# Volume cell as text (the sheet mixes thousand-separated values like
# "61.504.31" with plain decimals), so for those numeric-looking values
# we prefix with a leading apostrophe to force Excel to keep them as text
# instead of silently converting them to the Number type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.504.31"
$ws.Range("E2").Value = "  -2.44%  "
$ws.Range("D3").Value = "2.945.87"
$ws.Range("E3").Value = "  -3.52%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'582.75"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("D6").Value = "'140.30"
$ws.Range("E6").Value = "  -7.45%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -2.91%  "
$ws.Range("D9").Value = "2.941.06"
$ws.Range("E9").Value = "  -3.71%  "
$ws.Range("E10").Value = "  -6.41%  "
$ws.Range("D11").Value = "'5.69"
$ws.Range("E11").Value = "  -2.75%  "
$ws.Range("D12").Value = "'0.454"
$ws.Range("E12").Value = "  +1.28%  "
$ws.Range("E13").Value = "  -4.46%  "
$ws.Range("D14").Value = "'33.71"
$ws.Range("E14").Value = "  -6.65%  "
$ws.Range("E15").Value = "  +1.23%  "
$ws.Range("D16").Value = "3.436.46"
$ws.Range("E16").Value = "  -3.22%  "
$ws.Range("D17").Value = "'6.94"
$ws.Range("E17").Value = "  -2.84%  "
$ws.Range("D18").Value = "61.523.40"
$ws.Range("E18").Value = "  -2.35%  "
$ws.Range("D19").Value = "2.948.33"
$ws.Range("E19").Value = "  -3.41%  "
$ws.Range("D20").Value = "'446.46"
$ws.Range("E20").Value = "  -7.25%  "
$ws.Range("D21").Value = "'13.75"
$ws.Range("E21").Value = "  -3.89%  "
$ws.Range("D22").Value = "'0.676"
$ws.Range("E22").Value = "  -4.31%  "
$ws.Range("D23").Value = "'7.24"
$ws.Range("E23").Value = "  -3.46%  "
$ws.Range("D24").Value = "'80.63"
$ws.Range("E24").Value = "  -1.64%  "
$ws.Range("D25").Value = "'12.02"
$ws.Range("E25").Value = "  -5.09%  "
$ws.Range("D26").Value = "'2.13"
$ws.Range("E26").Value = "  -11.15%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").Value = "'9.53"
$ws.Range("E28").Value = "  -9.72%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("E30").Value = "  -2.03%  "
$ws.Range("D31").Value = "'6.82"
$ws.Range("E31").Value = "  -7.51%  "
$ws.Range("D32").Value = "'2.04"
$ws.Range("E32").Value = "  -7.33%  "
$ws.Range("D33").Value = "'26.88"
$ws.Range("E33").Value = "  -3.12%  "
$ws.Range("E34").Value = "  -4.29%  "
$ws.Range("E35").Value = "  -5.86%  "
$ws.Range("D36").Value = "0.0₃0763"
$ws.Range("E36").Value = "  -6.19%  "
$ws.Range("D37").Value = "'5.63"
$ws.Range("E37").Value = "  -4.81%  "
$ws.Range("E38").Value = "  -6.76%  "
$ws.Range("D39").Value = "'49.93"
$ws.Range("E39").Value = "  -0.98%  "
$ws.Range("D40").Value = "'9.07"
$ws.Range("E40").Value = "  -1.82%  "
$ws.Range("E41").Value = "  +2.42%  "
$ws.Range("D42").Value = "'2.76"
$ws.Range("E42").Value = "  -14.52%  "
$ws.Range("D43").Value = "'384.98"
$ws.Range("E43").Value = "  -10.24%  "
$ws.Range("E44").Value = "  -3.39%  "
$ws.Range("D45").Value = "2.702.45"
$ws.Range("E45").Value = "  -4.88%  "
$ws.Range("E46").Value = "  -8.94%  "
$ws.Range("D47").Value = "'36.63"
$ws.Range("E47").Value = "  -3.79%  "
$ws.Range("D48").Value = "'129.99"
$ws.Range("E48").Value = "  +2.15%  "
$ws.Range("E50").Value = "  -2.23%  "
$ws.Range("D51").Value = "'2.14"
$ws.Range("E51").Value = "  -2.74%  "
